$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values replacing the old "Strike#" values in column G
$kValues = @{
    2  = 2
    3  = 4
    4  = 3
    5  = 2
    6  = 4
    7  = 2
    8  = 1
    9  = 6
    10 = 4
    11 = 3
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 3
    18 = 2
    19 = 2
    20 = 3
    21 = 5
    22 = 0
    23 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
